$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the transaction date value in R2
$ws.Range("R2").Value = 20240907

# Update the view: scroll so column K is the top-left visible column, and
# move the selection to R3
$ws.Range("R3").Select()
$ws.Application.ActiveWindow.ScrollColumn = 11
